$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- C2: value + style (reuse format of F2, the CPF/verba style with yellow fill) ---
$ws.Range("F2").Copy()
$ws.Range("C2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C2").Value = 45128745242

# --- G2:G4: normalize date-column formatting (drop bold/courier font + wrap), keep DD/MM/YYYY + yellow fill ---
$ws.Range("A2").Copy()
$ws.Range("G2:G4").PasteSpecial(-4122)   # xlPasteFormats (font0/fill0/general base)
$ws.Range("G2:G4").Interior.Color = $ws.Range("F2").Interior.Color
$ws.Range("G2:G4").NumberFormat = "DD/MM/YYYY"
$ws.Range("G2:G4").WrapText = $false

$ws.Range("G2").Value = 43101
$ws.Range("G3").Value = 43101
$ws.Range("G4").Value = 43101

# --- View: scroll back to A1, move selection/cursor to B11 ---
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B11").Select()
